$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/5/2024  Through  2/11/2024"

# --- Structural cells: column layout shifted for rows 15, 26, 27 (style + type changes) ---
# Copy number-format (General/Number/Percent) from stable unaffected source cells, then set values.
# s=14 (General/text) source: D14   s=15 (#,##0) source: C36   s=16 (#,##0.0) source: K36
$ws.Range("D14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C15").PasteSpecial(-4163)

$ws.Range("C36").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1

$ws.Range("K36").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100

$ws.Range("C36").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1

$ws.Range("K36").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = 100

$ws.Range("D14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C26").PasteSpecial(-4163)

$ws.Range("C36").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 2

$ws.Range("K36").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 50

$excel.CutCopyMode = 0

# --- Remaining value-only updates ---
# Row 14
$ws.Range("N14").Value = -83.333333333333

# Row 15
$ws.Range("F15").Value = 2
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 66.666666666666
$ws.Range("L15").Value = 25
$ws.Range("M15").Value = 0

# Row 16
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 85.714285714285
$ws.Range("F16").Value = 39
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = 44.444444444444
$ws.Range("I16").Value = 57
$ws.Range("J16").Value = 36
$ws.Range("K16").Value = 58.333333333333
$ws.Range("L16").Value = 50
$ws.Range("M16").Value = -1.724137931034
$ws.Range("N16").Value = -74.666666666666

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 37
$ws.Range("G17").Value = 43
$ws.Range("H17").Value = -13.953488372093
$ws.Range("I17").Value = 67
$ws.Range("J17").Value = 68
$ws.Range("K17").Value = -1.470588235294
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 139.285714285714
$ws.Range("N17").Value = 63.414634146341

# Row 18
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -28.571428571428
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = 17.391304347826
$ws.Range("I18").Value = 37
$ws.Range("J18").Value = 30
$ws.Range("K18").Value = 23.333333333333
$ws.Range("L18").Value = 184.615384615385
$ws.Range("M18").Value = -19.565217391304
$ws.Range("N18").Value = -88.141025641025

# Row 19
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 100
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = 85.185185185185
$ws.Range("I19").Value = 129
$ws.Range("J19").Value = 86
$ws.Range("K19").Value = 50
$ws.Range("L19").Value = -41.891891891891
$ws.Range("M19").Value = 84.285714285714
$ws.Range("N19").Value = -11.643835616438

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 21
$ws.Range("H20").Value = -4.545454545454
$ws.Range("I20").Value = 32
$ws.Range("J20").Value = 34
$ws.Range("K20").Value = -5.882352941176
$ws.Range("L20").Value = 39.130434782608
$ws.Range("M20").Value = 77.777777777777
$ws.Range("N20").Value = -88.148148148148

# Row 21
$ws.Range("C21").Value = 48
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = 14.285714285714
$ws.Range("F21").Value = 226
$ws.Range("G21").Value = 171
$ws.Range("H21").Value = 32.163742690058
$ws.Range("I21").Value = 328
$ws.Range("J21").Value = 258
$ws.Range("K21").Value = 27.131782945736
$ws.Range("L21").Value = -10.869565217391
$ws.Range("M21").Value = 45.777777777777
$ws.Range("N21").Value = -67.232767232767

# Row 22
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 66.666666666666
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = 133.333333333333
$ws.Range("M22").Value = 16.666666666666

# Row 24
$ws.Range("C24").Value = 67
$ws.Range("D24").Value = 78
$ws.Range("E24").Value = -14.102564102564
$ws.Range("G24").Value = 246
$ws.Range("H24").Value = 18.292682926829
$ws.Range("I24").Value = 397
$ws.Range("J24").Value = 382
$ws.Range("K24").Value = 3.92670157068
$ws.Range("L24").Value = 48.134328358209
$ws.Range("M24").Value = 143.558282208589

# Row 25
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 58.333333333333
$ws.Range("F25").Value = 86
$ws.Range("G25").Value = 69
$ws.Range("H25").Value = 24.63768115942
$ws.Range("I25").Value = 122
$ws.Range("J25").Value = 103
$ws.Range("K25").Value = 18.446601941747
$ws.Range("L25").Value = 67.123287671232
$ws.Range("M25").Value = 121.818181818182

# Row 26
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 4
$ws.Range("K26").Value = 25
$ws.Range("L26").Value = 0

# Row 27
$ws.Range("C27").Value = 3
$ws.Range("F27").Value = 12
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 15
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 66.666666666666
